# Complete restructure and rewrite of documentation ready for v2
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet 1: "Notes" — update the specific-issue note text, and the
# active tab moves away from this sheet (handled by Activate() below
# on the "studies" sheet).
# ------------------------------------------------------------------
$wsNotes = $wb.Worksheets.Item("Notes")
$wsNotes.Range("A3").Value = "Specific issue: individual variant_num exceeds total_num"

# ------------------------------------------------------------------
# Sheet 2: "studies" — re-labelled header row, new example row, one
# extra column (G). Leave F2 (and its hyperlink) completely alone:
# its text/url is unchanged, only its shared-string slot moves.
# ------------------------------------------------------------------
$wsStudies = $wb.Worksheets.Item("studies")

$wsStudies.Range("B1").Value = "study_label"
$wsStudies.Range("C1").Value = "description"
$wsStudies.Range("D1").Value = "access_level"
$wsStudies.Range("E1").Value = "contributors"
$wsStudies.Range("F1").Value = "reference"
$wsStudies.Range("G1").Value = "reference_year"

$wsStudies.Range("A2").Value = "foo"
$wsStudies.Range("C2").ClearContents()
$wsStudies.Range("D2").Value = "public"

$wsStudies.Range("D3").Select()

# ------------------------------------------------------------------
# Sheet 3: "surveys" — full rebuild: two extra columns, every header
# loses the special font/number-format styling except the date-ish
# columns I:L which keep the "stored as text" number format.
# ------------------------------------------------------------------
$wsSurveys = $wb.Worksheets.Item("surveys")
$wsSurveys.Cells.Clear()

$wsSurveys.Range("A1").Value = "study_id"
$wsSurveys.Range("B1").Value = "survey_id"
$wsSurveys.Range("C1").Value = "country_name"
$wsSurveys.Range("D1").Value = "site_name"
$wsSurveys.Range("E1").Value = "latitude"
$wsSurveys.Range("F1").Value = "longitude"
$wsSurveys.Range("G1").Value = "location_method"
$wsSurveys.Range("H1").Value = "location_notes"
$wsSurveys.Range("I1").Value = "collection_start"
$wsSurveys.Range("J1").Value = "collection_end"
$wsSurveys.Range("K1").Value = "collection_day"
$wsSurveys.Range("L1").Value = "time_method"
$wsSurveys.Range("M1").Value = "time_notes"

$wsSurveys.Range("A1:H1").Style = "Normal"
$wsSurveys.Range("I1:L1").NumberFormat = "@"
$wsSurveys.Range("M1").Style = "Normal"

$wsSurveys.Range("A2").Value = "foo"
$wsSurveys.Range("B2").Value = "S01"
$wsSurveys.Range("E2").Value = 0
$wsSurveys.Range("F2").Value = 0
$wsSurveys.Range("H2").Value = "example data"
$wsSurveys.Range("K2").Value = "2020-01-01"
$wsSurveys.Range("M2").Value = "example data"

$wsSurveys.Range("A2:H2").Style = "Normal"
$wsSurveys.Range("K2").NumberFormat = "@"
$wsSurveys.Range("L2").NumberFormat = "@"
$wsSurveys.Range("M2").Style = "Normal"

$wsSurveys.Range("A1:M2").Select()

# ------------------------------------------------------------------
# Sheet 4: "counts" — same values, relabelled id columns.
# ------------------------------------------------------------------
$wsCounts = $wb.Worksheets.Item("counts")
$wsCounts.Range("A1").Value = "study_id"
$wsCounts.Range("B1").Value = "survey_id"
$wsCounts.Range("A2").Value = "foo"

$wsCounts.Range("B2").Select()

# ------------------------------------------------------------------
# Active sheet becomes "studies" (index 2, 0-based activeTab = 1) and
# it must carry tabSelected — Activate() drives both bookViews and
# the per-sheet sheetView flag in one go.
# ------------------------------------------------------------------
$wsStudies.Activate()
